$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header labels: "<name>_old" -> "<name>_FV2210" and
#    "<name>_new" -> "<name>_FV2304" (columns A1:U1 hold these headers,
#    also used as the Table1 column names below).
$headerRange = $ws.Range("A1:U66")
$headerRange.Replace("_old", "_FV2210", 2, 1, $false, $false, $false, $false)
$headerRange.Replace("_new", "_FV2304", 2, 1, $false, $false, $false, $false)

# 2. Turn the data range into an Excel Table ("Table1") with autofilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split below row 1, keep it in view while
#    scrolling) and select the first cell of the scrollable area.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
